# Apply "sourcePlot" column (J) to the "selections" sheet, and rework the
# I-column ("seedlot") formula for the headrow/pea rows (74-90) so that it
# builds the seedlot string from the new sourcePlot (J) column instead of a
# hard-coded "-1" suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("selections")

# --- New header for column J -------------------------------------------------
$ws.Range("J1").Value = "sourcePlot"

# --- J2:J73 - plot numbers already embedded in the existing I formulas ------
$sourcePlots = @{
    2 = 338;  3 = 249;  4 = 31;   5 = 97;   6 = 16;   7 = 174;  8 = 24;
    9 = 223;  10 = 56;  11 = 212; 12 = 128; 13 = 180; 14 = 10;  15 = 282;
    16 = 49;  17 = 367; 18 = 159; 19 = 100; 20 = 48;  21 = 247; 22 = 237;
    23 = 233; 24 = 264; 25 = 186; 26 = 251; 27 = 102; 28 = 161; 29 = 192;
    30 = 244; 31 = 316; 32 = 284; 33 = 341; 34 = 241; 35 = 253; 36 = 1;
    37 = 36;  38 = 153; 39 = 330; 40 = 3;   41 = 261; 42 = 95;  43 = 65;
    44 = 305; 45 = 134; 46 = 28;  47 = 194; 48 = 378; 49 = 203; 50 = 276;
    51 = 216; 52 = 182; 53 = 26;  54 = 297; 55 = 334; 56 = 227; 57 = 217;
    58 = 258; 59 = 323; 60 = 327; 61 = 37;  62 = 183; 63 = 199; 64 = 201;
    65 = 155; 66 = 347; 67 = 72;  68 = 313; 69 = 346; 70 = 349; 71 = 93;
    72 = 85;  73 = 370;
    74 = 38;  75 = 42;  76 = 52;  77 = 66;  78 = 68;  79 = 69;  80 = 70;
    81 = 71;  82 = 72;  83 = 74;  84 = 75;  85 = 76;  86 = 8;   87 = 83;
    88 = 93;  89 = 96
}

foreach ($row in $sourcePlots.Keys) {
    $ws.Cells.Item($row, 10).Value = $sourcePlots[$row]
}

# --- I column: seedlot formula --------------------------------------------
# Rows 2-73 (the "plot" rows) already build "<germplasm>-<obsUnit>" - untouched.
# Rows 74-90 (headrow / pea rows) move from the old
#   CONCAT(LEFT(C,38),"-",D,"-1")
# formula to
#   CONCAT(D,"-",LEFT(C,38),"-PLOT_",J)
# using the new sourcePlot column.
$ws.Range("I74:I90").Formula = '=_xlfn.CONCAT(D74,"-",LEFT(C74,38),"-PLOT_",J74)'

# --- sheet view state (matches the scroll/selection captured in the diff) ---
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("K89").Select()
